# Apply the edits described by the diff to the active workbook / sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Russian title in B1 loses the period right after "6.4.2.1"
$ws.Range("B1").Value = "6.4.2.1 Общий объем забора пресной воды "

# 2. Update the active cell / selection shown in the sheet view
$ws.Range("O2").Select()

# 3. Numeric value tweaks in column L (2022)
$ws.Range("L5").Value = 8741.9

# L7 used to hold the formula =L5-L8; it is now a plain static value
$ws.Range("L7").Value = 8483.5

$ws.Range("L14").Value = 1327.6

$ws.Range("L18").Value = 54
